$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header G1: "Parcel (ton)" -> "Revenue (Million)"
$ws.Range("G1").Value = "Revenue (Million)"

# Update values in column G (rows 7,8,9,12,14)
$ws.Range("G7").Value = 3.5
$ws.Range("G8").Value = 9
$ws.Range("G9").Value = 25
$ws.Range("G12").Value = 20
$ws.Range("G14").Value = 4

# Update the active cell / selection on the sheet view
$ws.Range("G19").Select()
